$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selected cell (sheet view selection)
$ws.Range("G19").Select()

# Update the input values (row 2)
$ws.Range("B2").Value = 50491
$ws.Range("F2").Value = 29079
$ws.Range("J2").Value = 40830
$ws.Range("N2").Value = 22925
$ws.Range("R2").Value = 53363

# Update the input values (row 3)
$ws.Range("B3").Value = 4922
$ws.Range("F3").Value = 17960
$ws.Range("J3").Value = 13625
$ws.Range("N3").Value = 36554
$ws.Range("R3").Value = 6846

$wb.Save()
